# NumberFormat.xlsx example update:
# Insert two new rows before the trailing "Some text" row, and populate them
# with a fraction-formatted numeric value (1.25), demonstrating the
# built-in "# ?/?" fraction format and a custom "# ?/100" fraction format.
# The pre-existing "Some text" row is pushed down by the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the last row ("Some text") down by inserting two blank rows above it.
$ws.Rows("19:20").Insert()

# Row 19: value 1.25 shown with the built-in fraction format "# ?/?" (numFmtId 12).
$ws.Cells.Item(19, 1).NumberFormat = "#\ ?/?"
$ws.Cells.Item(19, 1).Value = 1.25

# Row 20: value 1.25 shown with a custom fraction format "# ?/100".
$ws.Cells.Item(20, 1).NumberFormat = "#\ ?/100"
$ws.Cells.Item(20, 1).Value = 1.25

# Match the saved selection (active cell A19).
$ws.Range("A19").Select()
